# "dads matza cake image"
#
# 1. Add the missing recipe image for "עוגת מצות של פעם" (matza cake),
#    row 11 (id 10), into the image column (C).
# 2. Fix a typo in the קיצ'רי (kitchari) recipe's ingredients list:
#    "אדשים" -> "עדשים" (lentils), row 4 (id 3), ingredients column (E).
# 3. Leave the active cell selection on H4, matching the saved workbook view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C11").Value = "/recipes/images/8.jpg"
$ws.Range("E4").Value = "2 שיני שום|חמאה|&frac12;1 כוס אורז|1 כוס עדשים אדומות|3 כוסות מים|1 רסק עגבניות|פלפל|מלח"

$ws.Range("H4").Select()
